$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("OrdersPlaced")
$ws.Range("D1").Value = "FullPrice"
$ws.Range("F1").Value = "FinalPrice"
[void]$ws.Range("F1").Select()
